$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps being stored as text, matching the
# source data (values like "30.045.89" or "1.004" would otherwise be
# auto-coerced into numbers by Excel).

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.045.89"
$ws.Range("E2").Value = "  +0.69%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.884.88"
$ws.Range("E3").Value = "  -0.13%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("B5").Value = "XRP"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7355"
$ws.Range("E5").Value = "  -2.25%  "

$ws.Range("B6").Value = "BNB"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.48"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.23%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3163"
$ws.Range("E8").Value = "  +1.21%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07169"
$ws.Range("E9").Value = "  +0.74%  "

$ws.Range("B10").Value = "Solana"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.64"
$ws.Range("E10").Value = "  -2.67%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08343"
$ws.Range("E11").Value = "  -1.65%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.983.19"
$ws.Range("E12").Value = "  +5.16%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7547"
$ws.Range("E13").Value = "  -0.69%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.394"
$ws.Range("E14").Value = "  +0.55%  "

$ws.Range("B15").Value = "Litecoin"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.56"
$ws.Range("E15").Value = "  -0.86%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.076.97"
$ws.Range("E16").Value = "  +0.66%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.143"
$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "248.93"
$ws.Range("E18").Value = "  +2.31%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.54"
$ws.Range("E19").Value = "  -1.27%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007860"
$ws.Range("E20").Value = "  +0.66%  "

$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.26%  "

$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.144.97"
$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("E23").Value = "  +0.38%  "

$ws.Range("B24").Value = "Chainlink"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.906"
$ws.Range("E24").Value = "  -1.30%  "

$ws.Range("B25").Value = "Stellar"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1562"
$ws.Range("E25").Value = "  -1.96%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.271"
$ws.Range("E26").Value = "  -1.07%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.31"
$ws.Range("E27").Value = "  +0.98%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.63"
$ws.Range("E28").Value = "  -0.50%  "

$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.044"
$ws.Range("E29").Value = "  +0.63%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.478"
$ws.Range("E30").Value = "  -0.44%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.550"
$ws.Range("E31").Value = "  +1.02%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.535"
$ws.Range("E32").Value = "  +0.26%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.180"
$ws.Range("E33").Value = "  +0.47%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05323"
$ws.Range("E34").Value = "  -1.92%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.248"
$ws.Range("E35").Value = "  +0.51%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7666"
$ws.Range("E36").Value = "  +2.00%  "

$ws.Range("B37").Value = "Frax"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9960"
$ws.Range("E37").Value = "  -0.76%  "

$ws.Range("B38").Value = "HuobiToken"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.730"
$ws.Range("E38").Value = "  +0.68%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01953"
$ws.Range("E39").Value = "  +0.44%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.760"
$ws.Range("E40").Value = "  -0.44%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4546"
$ws.Range("E41").Value = "  +1.90%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.101.60"
$ws.Range("E42").Value = "  +0.44%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.046"
$ws.Range("E43").Value = "  -1.00%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.35"
$ws.Range("E44").Value = "  -0.38%  "

$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8747"
$ws.Range("E45").Value = "  +1.62%  "

$ws.Range("B46").Value = "PaxDollar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.004"
$ws.Range("E46").Value = "  +0.38%  "

$ws.Range("B47").Value = "Quant"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "104.22"
$ws.Range("E47").Value = "  +1.73%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.852"
$ws.Range("E48").Value = "  -0.45%  "

$ws.Range("B49").Value = "Aptos"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.554"
$ws.Range("E49").Value = "  -2.13%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.562"
$ws.Range("E50").Value = "  -1.99%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.046.82"
$ws.Range("E51").Value = "  +0.43%  "
